$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had 3 data rows:
#   Row1: Username | Password               (headers)
#   Row2: pankaj_kalra@unifyedqa.edu | Admin@2008s   (with hyperlinks)
#   Row3: philip_parker@unifyedqa.edu | Admin@2008   (with hyperlinks)
#
# The edit removes the "pankaj_kalra" row (old row 2), so the
# "philip_parker" row moves up to become row 2 and the sheet shrinks to
# A1:B2.

# 1) Move row 3's values into row 2 (overwriting the row to be removed).
$ws.Range("A2").Value2 = "philip_parker@unifyedqa.edu"
$ws.Range("B2").Value2 = "Admin@2008"

# 2) Remove the now-duplicate row 3 entirely, shrinking the used range.
$ws.Rows(3).Delete()

# 3) Re-point the hyperlinks so they match the data that is now in row 2
#    (clear the stale links first, then recreate the two that remain).
$ws.Hyperlinks.Delete()
[void]$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:philip_parker@unifyedqa.edu")
[void]$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Admin@2008")

# 4) Update the active selection left behind in the sheet view.
[void]$ws.Range("D11").Select()
